$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.538445
$ws.Range("H2").Value = 1.615335
$ws.Range("I2").Value = 0.03371608002174246
$ws.Range("J2").Value = 0.03371608002174246
$ws.Range("M2").Value = 133.3951123333333
$ws.Range("N2").Value = 400.185337
$ws.Range("O2").Value = 0.8984588679103155
$ws.Range("P2").Value = 0.8984588679103156
$ws.Range("Q2").Value = 71.82593126032167
$ws.Range("R2").Value = 646.433381342895
$ws.Range("S2").Value = 0.03029251108670834
$ws.Range("T2").Value = 0.03029251108670834
$ws.Range("G3").Value = 0.538445
$ws.Range("H3").Value = 1.615335
$ws.Range("I3").Value = 0.03371608002174246
$ws.Range("J3").Value = 0.03371608002174246
$ws.Range("M3").Value = 2.340788333333334
$ws.Range("N3").Value = 7.022365000000001
$ws.Range("O3").Value = 0.01576596023045448
$ws.Range("P3").Value = 0.01576596023045448
$ws.Range("Q3").Value = 1.260385774141667
$ws.Range("R3").Value = 11.343471967275
$ws.Range("S3").Value = 0.0005315663767496124
$ws.Range("T3").Value = 0.0005315663767496125
$ws.Range("G4").Value = 0.538445
$ws.Range("H4").Value = 1.615335
$ws.Range("I4").Value = 0.03371608002174246
$ws.Range("J4").Value = 0.03371608002174246
$ws.Range("M4").Value = 12.735128
$ws.Range("N4").Value = 38.205384
$ws.Range("O4").Value = 0.08577517185923002
$ws.Range("P4").Value = 0.08577517185923003
$ws.Range("Q4").Value = 6.85716599596
$ws.Range("R4").Value = 61.71449396364
$ws.Range("S4").Value = 0.002892002558284512
$ws.Range("T4").Value = 0.002892002558284512
$ws.Range("I5").Value = 0.7539416098905094
$ws.Range("J5").Value = 0.7539416098905093
$ws.Range("M5").Value = 133.3951123333333
$ws.Range("N5").Value = 400.185337
$ws.Range("O5").Value = 0.8984588679103155
$ws.Range("P5").Value = 0.8984588679103156
$ws.Range("Q5").Value = 1606.1344679266
$ws.Range("R5").Value = 14455.2102113394
$ws.Range("S5").Value = 0.6773855252927078
$ws.Range("T5").Value = 0.6773855252927078
$ws.Range("I6").Value = 0.7539416098905094
$ws.Range("J6").Value = 0.7539416098905093
$ws.Range("M6").Value = 2.340788333333334
$ws.Range("N6").Value = 7.022365000000001
$ws.Range("O6").Value = 0.01576596023045448
$ws.Range("P6").Value = 0.01576596023045448
$ws.Range("Q6").Value = 28.1840972920539
$ws.Range("R6").Value = 253.6568756284851
$ws.Range("S6").Value = 0.01188661343761859
$ws.Range("T6").Value = 0.0118866134376186
$ws.Range("I7").Value = 0.7539416098905094
$ws.Range("J7").Value = 0.7539416098905093
$ws.Range("M7").Value = 12.735128
$ws.Range("N7").Value = 38.205384
$ws.Range("O7").Value = 0.08577517185923002
$ws.Range("P7").Value = 0.08577517185923003
$ws.Range("Q7").Value = 153.3364129799974
$ws.Range("R7").Value = 1380.027716819976
$ws.Range("S7").Value = 0.064669471160183
$ws.Range("T7").Value = 0.064669471160183
$ws.Range("G8").Value = 3.391101666666666
$ws.Range("H8").Value = 10.173305
$ws.Range("I8").Value = 0.2123423100877482
$ws.Range("J8").Value = 0.2123423100877481
$ws.Range("M8").Value = 133.3951123333333
$ws.Range("N8").Value = 400.185337
$ws.Range("O8").Value = 0.8984588679103155
$ws.Range("P8").Value = 0.8984588679103156
$ws.Range("Q8").Value = 452.3563877587538
$ws.Range("R8").Value = 4071.207489828785
$ws.Range("S8").Value = 0.1907808315308994
$ws.Range("T8").Value = 0.1907808315308994
$ws.Range("G9").Value = 3.391101666666666
$ws.Range("H9").Value = 10.173305
$ws.Range("I9").Value = 0.2123423100877482
$ws.Range("J9").Value = 0.2123423100877481
$ws.Range("M9").Value = 2.340788333333334
$ws.Range("N9").Value = 7.022365000000001
$ws.Range("O9").Value = 0.01576596023045448
$ws.Range("P9").Value = 0.01576596023045448
$ws.Range("Q9").Value = 7.937851218480556
$ws.Range("R9").Value = 71.440660966325
$ws.Range("S9").Value = 0.00334778041608627
$ws.Range("T9").Value = 0.003347780416086271
$ws.Range("G10").Value = 3.391101666666666
$ws.Range("H10").Value = 10.173305
$ws.Range("I10").Value = 0.2123423100877482
$ws.Range("J10").Value = 0.2123423100877481
$ws.Range("M10").Value = 12.735128
$ws.Range("N10").Value = 38.205384
$ws.Range("O10").Value = 0.08577517185923002
$ws.Range("P10").Value = 0.08577517185923003
$ws.Range("Q10").Value = 43.18611378601333
$ws.Range("R10").Value = 388.67502407412
$ws.Range("S10").Value = 0.01821369814076251
$ws.Range("T10").Value = 0.01821369814076251
